# Update the "取得日時" (acquired timestamp) column A for data rows 2-13
# on the "ランサーズ" sheet from "2025-12-25 12:50:27" to "2025-12-25 18:26:42".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-12-25 12:50:27"
$newValue = "2025-12-25 18:26:42"

for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
